$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B43").Value = 3.27
$ws.Range("C43").Value = 2.36
$ws.Range("D43").Value = 1.74
$ws.Range("E43").Value = 1.24
$ws.Range("F43").Value = 0.77

$ws.Range("B37").Select()
